# Added meta type for std::weak_ptr
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta types")

# Insert a new row above row 42 (existing rows 42..63 shift down to 43..64)
$ws.Rows.Item(42).Insert()

# Fill in the new row with the std::weak_ptr<T> meta-type entry
$ws.Cells.Item(42, 1).Value = "tkStdWeakPtr"
$ws.Cells.Item(42, 2).Value = 104
$ws.Cells.Item(42, 3).Value = "std::weak_ptr<T>"
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(42, 6).Value = "T"

# Update the view state (scroll position + selection) to match the authored file
$excel.Goto($ws.Range("A30"), $true)
[void]$ws.Range("B43").Select()
